$d = $word.ActiveDocument

# --- 1) Heading paragraph: wrap "Techerudite_Test" runs with spell-check
#        proofErr markers (spellStart / spellEnd) ---------------------------
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004500E6" w:rsidRDefault="004500E6" w:rsidP="004500E6"><w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="center"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Tech</w:t></w:r><w:r w:rsidR="00093446"><w:t>eru</w:t></w:r><w:r><w:t>dite_Test</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(1).Range.InsertXML($headingXml)

# --- 2) Closing paragraph: split "Thank You !" into "Thank " + "You ",
#        drop the "!" and wrap "You " with grammar-check proofErr markers --
$lastIndex = $d.Paragraphs.Count
$closingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00093446" w:rsidRPr="004500E6" w:rsidRDefault="00093446" w:rsidP="004500E6"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">There is also improvement in it but according to keep in deadline I have keep this version </w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Thank </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">You </w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$d.Paragraphs($lastIndex).Range.InsertXML($closingXml)
